$d = $word.ActiveDocument

# Replace field-name tokens inside the laboratory report merge fields.
$d.Content.Find.Execute("collectionDateTime", $true, $false, $false, $false, $false,
                         $true, 1, $false, "date", 2)
$d.Content.Find.Execute("batteryType", $true, $false, $false, $false, $false,
                         $true, 1, $false, "test", 2)
$d.Content.Find.Execute("testStatus", $true, $false, $false, $false, $false,
                         $true, 1, $false, "status", 2)
